$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ทำการบ้านเลขไม่ได้เลยทำยังไงดี"
$ws.Range("B2").Value = "ทำการบ้าน"
$ws.Range("C2").Value = "การเรียน"
$ws.Range("D2").Value = "เลข"
$ws.Range("E2").Value = "คณิตศาสตร์"
$ws.Range("F2").Value = "การบ้าน"
$ws.Range("G2").ClearContents()

# Row 3
$ws.Range("A3").Value = "คณิตศาสตร์ยากจัง"
$ws.Range("B3").Value = "คณิตศาสตร์"
$ws.Range("C3").Value = "การบ้าน"
$ws.Range("D3:H3").ClearContents()

# Row 4 (new)
$ws.Range("A4").Value = "ทำไมการเรียนออนไลน์ถึงมีการบ้านเยอะจัง"
$ws.Range("B4").Value = "การเรียน"
$ws.Range("C4").Value = "การบ้าน"

# Row 5 (new)
$ws.Range("A5").Value = "ผมเกลียดเลข แต่แม่บังคับให้เรียนพิเศษเพิ่ม ทำอย่างไรดี"
$ws.Range("B5").Value = "เลข"
$ws.Range("C5").Value = "คณิตศาสตร์"
$ws.Range("D5").Value = "แม่"
$ws.Range("E5").Value = "ครอบครัว"
